# Update automàtic: dades i banners [2026-02-25 20:49]
# Refreshes the DATA_EXTRACCIO timestamps and the observation values that
# meteo.cat reported slightly differently on the 20:48 re-poll.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = "2026-02-25 20:48:17"
# Row 3
$ws.Range("E3").Value2 = "2026-02-25 20:48:19"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value2 = "37%"
# Row 4
$ws.Range("E4").Value2 = "2026-02-25 20:48:22"
$ws.Range("O4").Value2 = "8.7 °C"
# Row 5
$ws.Range("E5").Value2 = "2026-02-25 20:48:24"
$ws.Range("N5").Value2 = "2.3 °C 20:01 TU"
# Row 6
$ws.Range("E6").Value2 = "2026-02-25 20:48:26"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value2 = "90%"
$ws.Range("J6").Value2 = "1021.7 hPa"
# Row 7
$ws.Range("E7").Value2 = "2026-02-25 20:48:29"
# Row 8
$ws.Range("E8").Value2 = "2026-02-25 20:48:31"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value2 = "79%"
$ws.Range("J8").Value2 = "1021.0 hPa"
$ws.Range("O8").Value2 = "11.9 °C"
# Row 9
$ws.Range("E9").Value2 = "2026-02-25 20:48:34"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value2 = "91%"
$ws.Range("O9").Value2 = "10.2 °C"
# Row 10
$ws.Range("E10").Value2 = "2026-02-25 20:48:35"
# Row 11
$ws.Range("E11").Value2 = "2026-02-25 20:48:36"
$ws.Range("O11").Value2 = "9.0 °C"
# Row 12
$ws.Range("E12").Value2 = "2026-02-25 20:48:37"
$ws.Range("O12").Value2 = "10.2 °C"
# Row 13
$ws.Range("E13").Value2 = "2026-02-25 20:48:38"
$ws.Range("J13").Value2 = "1022.7 hPa"
# Row 14
$ws.Range("E14").Value2 = "2026-02-25 20:48:39"
$ws.Range("O14").Value2 = "10.6 °C"
# Row 15
$ws.Range("E15").Value2 = "2026-02-25 20:48:40"
# Row 16
$ws.Range("E16").Value2 = "2026-02-25 20:48:41"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value2 = "31%"
$ws.Range("O16").Value2 = "3.2 °C"
# Row 17
$ws.Range("E17").Value2 = "2026-02-25 20:48:42"
$ws.Range("N17").Value2 = "5.5 °C 20:10 TU"
$ws.Range("O17").Value2 = "9.3 °C"
# Row 18
$ws.Range("E18").Value2 = "2026-02-25 20:48:44"
# Row 19
$ws.Range("E19").Value2 = "2026-02-25 20:48:45"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value2 = "49%"
$ws.Range("N19").Value2 = "8.3 °C 20:29 TU"
$ws.Range("O19").Value2 = "12.3 °C"
# Row 20
$ws.Range("E20").Value2 = "2026-02-25 20:48:46"
$ws.Range("N20").Value2 = "-0.3 °C 20:25 TU"
$ws.Range("O20").Value2 = "2.8 °C"
# Row 21
$ws.Range("E21").Value2 = "2026-02-25 20:48:48"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value2 = "57%"
$ws.Range("K21").Value2 = "15.7 MJ/m2"
# Row 22
$ws.Range("E22").Value2 = "2026-02-25 20:48:50"
$ws.Range("N22").Value2 = "0.1 °C 20:26 TU"
$ws.Range("O22").Value2 = "2.5 °C"
# Row 23
$ws.Range("E23").Value2 = "2026-02-25 20:48:53"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value2 = "34%"
$ws.Range("N23").Value2 = "1.8 °C 20:29 TU"
$ws.Range("O23").Value2 = "4.0 °C"
# Row 24
$ws.Range("E24").Value2 = "2026-02-25 20:48:55"
$ws.Range("J24").Value2 = "1020.0 hPa"
# Row 25
$ws.Range("E25").Value2 = "2026-02-25 20:48:58"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value2 = "34%"
$ws.Range("O25").Value2 = "5.3 °C"
# Row 26
$ws.Range("E26").Value2 = "2026-02-25 20:49:00"
$ws.Range("J26").Value2 = "1019.5 hPa"
$ws.Range("N26").Value2 = "5.3 °C 20:16 TU"
$ws.Range("O26").Value2 = "10.1 °C"
# Row 27
$ws.Range("E27").Value2 = "2026-02-25 20:49:03"
$ws.Range("O27").Value2 = "5.3 °C"
# Row 28
$ws.Range("E28").Value2 = "2026-02-25 20:49:05"
$ws.Range("J28").Value2 = "1021.8 hPa"
# Row 29
$ws.Range("E29").Value2 = "2026-02-25 20:49:08"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value2 = "91%"
# Row 30
$ws.Range("E30").Value2 = "2026-02-25 20:49:10"
$ws.Range("J30").Value2 = "1021.8 hPa"
$ws.Range("O30").Value2 = "10.6 °C"
# Row 31
$ws.Range("E31").Value2 = "2026-02-25 20:49:12"
# Row 32
$ws.Range("E32").Value2 = "2026-02-25 20:49:15"
$ws.Range("O32").Value2 = "9.4 °C"
# Row 33
$ws.Range("E33").Value2 = "2026-02-25 20:49:17"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value2 = "51%"
# Row 34
$ws.Range("E34").Value2 = "2026-02-25 20:49:20"
$ws.Range("N34").Value2 = "-0.5 °C 20:01 TU"
$ws.Range("O34").Value2 = "3.5 °C"
# Row 35
$ws.Range("E35").Value2 = "2026-02-25 20:49:22"
$ws.Range("J35").Value2 = "1019.5 hPa"
$ws.Range("O35").Value2 = "12.5 °C"
# Row 36
$ws.Range("E36").Value2 = "2026-02-25 20:49:24"
$ws.Range("J36").Value2 = "1021.9 hPa"
# Row 37
$ws.Range("E37").Value2 = "2026-02-25 20:49:27"
$ws.Range("J37").Value2 = "1023.4 hPa"
$ws.Range("O37").Value2 = "6.7 °C"
# Row 38
$ws.Range("E38").Value2 = "2026-02-25 20:49:29"
# Row 39
$ws.Range("E39").Value2 = "2026-02-25 20:49:32"
# Row 40
$ws.Range("E40").Value2 = "2026-02-25 20:49:34"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value2 = "60%"
$ws.Range("J40").Value2 = "1021.8 hPa"
$ws.Range("O40").Value2 = "9.7 °C"
# Row 41
$ws.Range("E41").Value2 = "2026-02-25 20:49:36"
$ws.Range("J41").Value2 = "1021.0 hPa"
$ws.Range("O41").Value2 = "11.6 °C"
# Row 42
$ws.Range("E42").Value2 = "2026-02-25 20:49:39"
# Row 43
$ws.Range("E43").Value2 = "2026-02-25 20:49:41"
# Row 44
$ws.Range("E44").Value2 = "2026-02-25 20:49:43"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value2 = "45%"
# Row 45
$ws.Range("E45").Value2 = "2026-02-25 20:49:46"
$ws.Range("J45").Value2 = "1019.8 hPa"
$ws.Range("O45").Value2 = "11.0 °C"
# Row 46
$ws.Range("E46").Value2 = "2026-02-25 20:49:48"
$ws.Range("O46").Value2 = "9.8 °C"
